# Updated remaining queries for C3DC
# The JOINs across df_study/df_participant/df_diagnoses/df_treatments/
# df_treatment_resp/df_survival/df_reference_files used to key off the
# generic "id" column; the upstream data model renamed the join keys to
# the explicit "study_id" / "participant_id" columns. Apply that rename
# to every query still using the old column names (StudiesTab, StatQuery,
# ParticipantsTab, DiagnosisTab, TreatmentTab, TreatmentRespTab,
# SurvivalTab).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$queryCells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellRef in $queryCells) {
    $rng = $ws.Range($cellRef)
    $sql = $rng.Text

    $sql = $sql.Replace('df_participant prt ON std.id = prt."study.id"', 'df_participant prt ON std.study_id = prt."study.study_id"')
    $sql = $sql.Replace('df_diagnoses dgn ON prt.id = dgn."participant.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"')
    $sql = $sql.Replace('df_treatments trt ON prt.id = trt."participant.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"')
    $sql = $sql.Replace('df_treatment_resp trr ON prt.id = trr."participant.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"')
    $sql = $sql.Replace('df_survival srv ON prt.id = srv."participant.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"')
    $sql = $sql.Replace('df_reference_files rfs ON std.id = rfs."study.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"')

    $rng.Value = $sql
}

# Widen column C (StatQuery) to fit the updated text, matching the
# author's manual resize; drop the stale bestFit autosize.
$ws.Columns.Item(3).ColumnWidth = 68.1666666666667

# Author left the selection on B2 (top of the query list) instead of the
# previous scrolled-down C4 selection.
$ws.Range("B2").Select()
